$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its original text (string) representation:
# several updated prices (e.g. "566.40", "7.20", "0.0000234") look like
# numbers and would otherwise be auto-converted/reformatted by Excel,
# losing trailing zeros or switching to scientific notation.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.845.76"
$ws.Range("E2").Value = "  +3.90%  "
$ws.Range("D3").Value = "3.016.74"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "566.40"
$ws.Range("E5").Value = "  +3.18%  "
$ws.Range("D6").Value = "141.21"
$ws.Range("E6").Value = "  +8.37%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +2.36%  "
$ws.Range("D9").Value = "3.009.10"
$ws.Range("E9").Value = "  +2.68%  "
$ws.Range("E10").Value = "  +6.65%  "
$ws.Range("D11").Value = "5.32"
$ws.Range("E11").Value = "  +11.80%  "
$ws.Range("E12").Value = "  +3.33%  "
$ws.Range("D13").Value = "0.0000234"
$ws.Range("E13").Value = "  +5.55%  "
$ws.Range("D14").Value = "34.14"
$ws.Range("E14").Value = "  +3.92%  "
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("D16").Value = "3.512.82"
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("D17").Value = "7.20"
$ws.Range("E17").Value = "  +4.60%  "
$ws.Range("D18").Value = "3.012.53"
$ws.Range("E18").Value = "  +2.69%  "
$ws.Range("D19").Value = "59.806.25"
$ws.Range("E19").Value = "  +3.78%  "
$ws.Range("D20").Value = "440.28"
$ws.Range("D21").Value = "13.69"
$ws.Range("E21").Value = "  +3.53%  "
$ws.Range("D22").Value = "0.723"
$ws.Range("E22").Value = "  +5.14%  "
$ws.Range("D23").Value = "7.14"
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("D24").Value = "13.38"
$ws.Range("E24").Value = "  +2.34%  "
$ws.Range("D25").Value = "80.88"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "2.24"
$ws.Range("E27").Value = "  +12.53%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("E29").Value = "  +3.70%  "
$ws.Range("D30").Value = "7.88"
$ws.Range("E30").Value = "  +5.69%  "
$ws.Range("D31").Value = "6.35"
$ws.Range("E31").Value = "  +5.72%  "
$ws.Range("D32").Value = "26.08"
$ws.Range("E32").Value = "  +3.62%  "
$ws.Range("D33").Value = "0.105"
$ws.Range("E33").Value = "  +8.14%  "
$ws.Range("D34").Value = "0.0₃0796"
$ws.Range("E34").Value = "  +16.39%  "
$ws.Range("E35").Value = "  +6.83%  "
$ws.Range("D36").Value = "5.95"
$ws.Range("E36").Value = "  +5.19%  "
$ws.Range("D37").Value = "2.12"
$ws.Range("E37").Value = "  +2.21%  "
$ws.Range("D38").Value = "49.34"
$ws.Range("E38").Value = "  +2.10%  "
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").Value = "2.83"
$ws.Range("E40").Value = "  +10.88%  "
$ws.Range("D41").Value = "406.01"
$ws.Range("E41").Value = "  +7.98%  "
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("D43").Value = "2.775.78"
$ws.Range("E43").Value = "  +2.98%  "
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("E45").Value = "  +6.76%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "123.02"
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("D48").Value = "2.05"
$ws.Range("E48").Value = "  +3.83%  "
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("D50").Value = "34.15"
$ws.Range("E50").Value = "  +20.83%  "
$ws.Range("D51").Value = "23.73"
